# Generate Report for Handoff
# b.md has now been handed off for both zh-cn and de-de locales.
# Update the Overview sheet and the per-locale sheets accordingly.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-24 16:38:49"

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-24 16:38:44"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d84d36e3566035b7f7e21a98e99f02b618c2fdd/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9c517f98aaf58df8b3d54ce11a134ff593f3a85/e2e/b.md."
$zhcn.Range("P1").ColumnWidth = 39.17

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-24 16:38:49"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d84d36e3566035b7f7e21a98e99f02b618c2fdd/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9c517f98aaf58df8b3d54ce11a134ff593f3a85/e2e/b.md."
$dede.Range("P1").ColumnWidth = 39.17
